$wb = $excel.ActiveWorkbook

# --- Sheet 2: Restricciones_del_lider ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").NumberFormat = "@"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("A2").Value = "1.0499999999999998 - x"
$ws.Range("B2").Value = "-2.05"
$ws.Range("D2").Value = "0.24"
$ws.Range("A3").Value = "-1.05 + x"
$ws.Range("B3").Value = "0.050000000000000044"
$ws.Range("D3").Value = "0.72"
$ws.Range("B2:B3").Style = "Normal"
$ws.Range("D2:D3").Style = "Normal"

# --- Sheet 3: Restricciones_del_follower ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("A3").NumberFormat = "@"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("A2").Value = "0"
$ws.Range("B2").Value = "-1"
$ws.Range("D2").Value = "0.47"
$ws.Range("E2").Value = "0"
$ws.Range("F2").Value = "0"
$ws.Range("A3").Value = "0"
$ws.Range("B3").Value = "-1"
$ws.Range("D3").Value = "0.88"
$ws.Range("E3").Value = "0"
$ws.Range("F3").Value = "0"
$ws.Range("A2:F3").Style = "Normal"

# --- Sheet 4: Punto_modificado ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("A2").Value = "1.05"
$ws.Range("B2").Value = "2.85"
$ws.Range("A2:B2").Style = "Normal"

# --- Sheet 5: Vector_bf ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "-1.834875"
$ws.Range("A2").Style = "Normal"

# --- Sheet 6: Vector_BF ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A2").Value = "-0.45062500000000005"
$ws.Range("A3").Value = "-1.3625000000000003"
$ws.Range("A2:A3").Style = "Normal"

# --- Sheet 7: Vector_Alpha ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("A2").Value = 2.25
